$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.076.59"
Set-TextValue "E2" "  +0.75%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.130.52"
Set-TextValue "E3" "  -0.18%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "530.64"
Set-TextValue "E5" "  +0.54%  "

# Row 6 - Solana
Set-TextValue "D6" "138.73"
Set-TextValue "E6" "  -0.76%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.08%  "

# Row 8 - XRP
Set-TextValue "D8" "0.461"
Set-TextValue "E8" "  +4.07%  "

# Row 9 - Toncoin
Set-TextValue "E9" "  +1.51%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  -0.97%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.408"
Set-TextValue "E11" "  +3.00%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "3.670.62"
Set-TextValue "E12" "  -0.45%  "

# Row 13 - TRON
Set-TextValue "E13" "  +1.21%  "

# Row 14 - Avalanche
Set-TextValue "D14" "25.55"
Set-TextValue "E14" "  +0.15%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000163"
Set-TextValue "E15" "  -1.04%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "58.147.17"
Set-TextValue "E16" "  +0.57%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.125.19"
Set-TextValue "E17" "  -0.78%  "

# Row 18 - Polkadot
Set-TextValue "D18" "6.00"
Set-TextValue "E18" "  -1.54%  "

# Row 19 - Chainlink
Set-TextValue "D19" "12.70"
Set-TextValue "E19" "  -0.94%  "

# Row 20 - Uniswap
Set-TextValue "D20" "8.08"
Set-TextValue "E20" "  +2.00%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "352.92"
Set-TextValue "E21" "  -0.52%  "

# Row 22 - Dai
Set-TextValue "D22" "1.00"
Set-TextValue "E22" "  +0.17%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.12"
Set-TextValue "E23" "  +0.54%  "

# Row 24 - Polygon
Set-TextValue "D24" "0.505"
Set-TextValue "E24" "  -0.71%  "

# Row 25 - Kaspa
Set-TextValue "D25" "0.168"
Set-TextValue "E25" "  -1.35%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.23%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0₃0886"
Set-TextValue "E27" "  -4.34%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "7.29"
Set-TextValue "E28" "  -1.64%  "

# Row 29 - RenderToken
Set-TextValue "D29" "6.11"
Set-TextValue "E29" "  -3.53%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  -1.41%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "21.31"
Set-TextValue "E31" "  +0.80%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "4.99"
Set-TextValue "E32" "  +1.66%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.16"
Set-TextValue "E33" "  -2.72%  "

# Row 34 - Monero
Set-TextValue "D34" "158.82"
Set-TextValue "E34" "  +0.75%  "

# Row 35 - Aptos
Set-TextValue "D35" "6.07"
Set-TextValue "E35" "  -1.74%  "

# Row 36 - EnergySwap
Set-TextValue "D36" "26.31"
Set-TextValue "E36" "  +1.10%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "1.26"
Set-TextValue "E37" "  -1.09%  "

# Row 38 - Stacks
Set-TextValue "D38" "1.68"
Set-TextValue "E38" "  +3.69%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.0671"
Set-TextValue "E39" "  +0.18%  "

# Row 40 - now Filecoin (was Mantle)
Set-TextValue "B40" "Filecoin"
Set-TextValue "C40" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D40" "4.00"
Set-TextValue "E40" "  -2.51%  "

# Row 41 - now Mantle (was Filecoin)
Set-TextValue "B41" "Mantle"
Set-TextValue "C41" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D41" "0.698"
Set-TextValue "E41" "  -0.81%  "

# Row 42 - OKB
Set-TextValue "D42" "37.56"
Set-TextValue "E42" "  +2.56%  "

# Row 43 - Maker
Set-TextValue "D43" "2.394.00"
Set-TextValue "E43" "  +2.79%  "

# Row 44 - RenzoRestakedETH
Set-TextValue "D44" "3.162.33"
Set-TextValue "E44" "  -0.72%  "

# Row 45 - FirstDigitalUSD
Set-TextValue "D45" "1.00"
Set-TextValue "E45" "  -0.02%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0267"
Set-TextValue "E46" "  -2.33%  "

# Row 47 - ONDO
Set-TextValue "D47" "0.978"
Set-TextValue "E47" "  -0.69%  "

# Row 48 - Cosmos
Set-TextValue "D48" "6.04"
Set-TextValue "E48" "  -0.32%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "19.83"
Set-TextValue "E49" "  -1.81%  "

# Row 50 - SuiNetwork
Set-TextValue "D50" "0.739"
Set-TextValue "E50" "  -1.76%  "

# Row 51 - Stellar
Set-TextValue "D51" "0.0908"
Set-TextValue "E51" "  +1.80%  "
